$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update 想去人数 (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 340
$ws1.Range("F4").Value = 2896
$ws1.Range("F5").Value = 71
$ws1.Range("F6").Value = 610

# Sheet "全部类型" (fourth sheet) - update 想去人数 (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 340
$ws4.Range("F6").Value = 2896
$ws4.Range("F7").Value = 71
$ws4.Range("F8").Value = 610
